# Auto-generated edit script for horarios workbook update
# Applies new scraped rows to LP1912 (159 rows), LP1912-215 (21 rows),
# and refreshes the "Ultima actualizacion" timestamp on 6203-6173.
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 10:52:37"
$ws1.Range("A3").Value = "Total filas: 159"

$data1 = @(
    @(6, '00:09:19', '01:12', '215_ALUAR', 63, 'LP1912'),
    @(7, '00:09:19', '01:58', '14_ABASTO', 109, 'LP1912'),
    @(8, '02:17:56', '02:57', '215_ALUAR', 40, 'LP1912'),
    @(9, '01:16:09', '02:58', '215_ALUAR', 102, 'LP1912'),
    @(10, '01:55:40', '03:48', '14_ABASTO', 113, 'LP1912'),
    @(11, '02:17:56', '04:01', '81_EL PELIGRO', 104, 'LP1912'),
    @(12, '04:41:47', '04:45', '215A_EL PATO', 4, 'LP1912'),
    @(13, '02:56:23', '04:46', '215A_EL PATO', 110, 'LP1912'),
    @(14, '02:56:23', '04:53', '11_ETCHEVERRY', 117, 'LP1912'),
    @(15, '04:54:25', '04:54', '11_ETCHEVERRY', 0, 'LP1912'),
    @(16, '03:25:20', '05:16', '17_ROMERO', 111, 'LP1912'),
    @(17, '03:25:20', '05:22', '23_HERNANDEZ', 117, 'LP1912'),
    @(18, '05:18:42', '05:25', '23_HERNANDEZ', 7, 'LP1912'),
    @(19, '04:41:47', '05:34', '215B_EL PATO', 53, 'LP1912'),
    @(20, '03:55:07', '05:35', '215B_EL PATO', 100, 'LP1912'),
    @(21, '03:55:07', '05:46', '15_ABASTO', 111, 'LP1912'),
    @(22, '03:55:07', '05:54', '10_OLMOS', 119, 'LP1912'),
    @(23, '04:19:04', '06:04', '16_SANTA ANA', 105, 'LP1912'),
    @(24, '04:41:47', '06:11', '215A_EL PATO', 90, 'LP1912'),
    @(25, '04:19:04', '06:12', '215A_EL PATO', 113, 'LP1912'),
    @(26, '04:19:04', '06:14', '225_HARAS DEL SUR', 115, 'LP1912'),
    @(27, '06:15:33', '06:15', '225_HARAS DEL SUR', 0, 'LP1912'),
    @(28, '04:41:47', '06:21', '26_HERNANDEZ', 100, 'LP1912'),
    @(29, '04:41:47', '06:27', '23_HERNANDEZ', 106, 'LP1912'),
    @(30, '04:41:47', '06:29', '86_EST CHICA-ESC AGRARIA', 108, 'LP1912'),
    @(31, '04:41:47', '06:31', '16_SANTA ANA', 110, 'LP1912'),
    @(32, '04:54:25', '06:44', '225_C ROCA-H SUR', 110, 'LP1912'),
    @(33, '04:54:25', '06:46', '215C_EL PATO', 112, 'LP1912'),
    @(34, '05:18:42', '06:59', '14_ABASTO', 101, 'LP1912'),
    @(35, '06:58:31', '07:00', '14_ABASTO', 2, 'LP1912'),
    @(36, '06:15:33', '07:01', '16_SANTA ANA', 46, 'LP1912'),
    @(37, '05:53:46', '07:04', '23_HERNANDEZ', 71, 'LP1912'),
    @(38, '05:18:42', '07:05', '15_ABASTO', 107, 'LP1912'),
    @(39, '06:44:40', '07:05', '23_HERNANDEZ', 21, 'LP1912'),
    @(40, '05:18:42', '07:07', '225_GOMEZ', 109, 'LP1912'),
    @(41, '06:58:31', '07:08', '15_ABASTO', 10, 'LP1912'),
    @(42, '06:44:40', '07:09', '15_ABASTO', 25, 'LP1912'),
    @(43, '05:18:42', '07:11', '215A_EL PATO', 113, 'LP1912'),
    @(44, '05:18:42', '07:15', '11_ETCHEVERRY', 117, 'LP1912'),
    @(45, '06:44:40', '07:16', '16_SANTA ANA', 32, 'LP1912'),
    @(46, '05:53:46', '07:21', '26_HERNANDEZ', 88, 'LP1912'),
    @(47, '06:15:33', '07:23', '10_OLMOS', 68, 'LP1912'),
    @(48, '05:53:46', '07:31', '11_ETCHEVERRY', 98, 'LP1912'),
    @(49, '07:24:45', '07:31', '16_SANTA ANA', 7, 'LP1912'),
    @(50, '05:53:46', '07:32', '84_COLONIA URQUIZA-ESC 49', 99, 'LP1912'),
    @(51, '07:24:45', '07:34', '23_HERNANDEZ', 10, 'LP1912'),
    @(52, '05:53:46', '07:36', '27_EL RETIRO', 103, 'LP1912'),
    @(53, '05:53:46', '07:39', '10_OLMOS', 106, 'LP1912'),
    @(54, '07:24:45', '07:46', '16_SANTA ANA', 22, 'LP1912'),
    @(55, '05:53:46', '07:47', '14_ABASTO', 114, 'LP1912'),
    @(56, '07:48:31', '07:48', '14_ABASTO', 0, 'LP1912'),
    @(57, '05:53:46', '07:51', '215D_EL PATO', 118, 'LP1912'),
    @(58, '06:58:31', '07:58', '16_SANTA ANA', 60, 'LP1912'),
    @(59, '07:24:45', '08:03', '11_ETCHEVERRY', 39, 'LP1912'),
    @(60, '06:15:33', '08:07', '16_SANTA ANA', 112, 'LP1912'),
    @(61, '07:48:31', '08:10', '16_SANTA ANA', 22, 'LP1912'),
    @(62, '06:15:33', '08:12', '15_ABASTO', 117, 'LP1912'),
    @(63, '06:58:31', '08:13', '10_OLMOS', 75, 'LP1912'),
    @(64, '06:44:40', '08:21', '26_HERNANDEZ', 97, 'LP1912'),
    @(65, '06:44:40', '08:22', '16_P MOR-SANTA ANA', 98, 'LP1912'),
    @(66, '06:44:40', '08:23', '215B_EL PATO', 99, 'LP1912'),
    @(67, '06:44:40', '08:27', '84_COLONIA URQUIZA-ESC 49', 103, 'LP1912'),
    @(68, '08:31:53', '08:31', '10_OLMOS', 0, 'LP1912'),
    @(69, '07:48:31', '08:32', '23_HERNANDEZ', 44, 'LP1912'),
    @(70, '08:00:32', '08:33', '10_OLMOS', 33, 'LP1912'),
    @(71, '06:58:31', '08:34', '23_HERNANDEZ', 96, 'LP1912'),
    @(72, '06:44:40', '08:35', '23_HERNANDEZ', 111, 'LP1912'),
    @(73, '06:44:40', '08:42', '81_EL PELIGRO', 118, 'LP1912'),
    @(74, '07:24:45', '08:43', '14_ABASTO', 79, 'LP1912'),
    @(75, '08:31:53', '08:46', '16_SANTA ANA', 15, 'LP1912'),
    @(76, '08:47:51', '08:47', '16_SANTA ANA', 0, 'LP1912'),
    @(77, '07:48:31', '08:53', '10_OLMOS', 65, 'LP1912'),
    @(78, '06:58:31', '08:54', '17_ROMERO', 116, 'LP1912'),
    @(79, '08:55:25', '08:55', '10_OLMOS', 0, 'LP1912'),
    @(80, '08:55:25', '08:56', '16_SANTA ANA', 1, 'LP1912'),
    @(81, '07:24:45', '09:01', '215A_EL PATO', 97, 'LP1912'),
    @(82, '08:00:32', '09:03', '11_ETCHEVERRY', 63, 'LP1912'),
    @(83, '08:31:53', '09:04', '23_HERNANDEZ', 33, 'LP1912'),
    @(84, '08:47:51', '09:05', '23_HERNANDEZ', 18, 'LP1912'),
    @(85, '08:55:25', '09:06', '23_HERNANDEZ', 11, 'LP1912'),
    @(86, '07:48:31', '09:07', '23_HERNANDEZ', 79, 'LP1912'),
    @(87, '08:00:32', '09:08', '23_HERNANDEZ', 68, 'LP1912'),
    @(88, '07:24:45', '09:10', '16_P MOR-SANTA ANA', 106, 'LP1912'),
    @(89, '08:47:51', '09:13', '16_SANTA ANA', 26, 'LP1912'),
    @(90, '07:24:45', '09:16', '27_EL RETIRO', 112, 'LP1912'),
    @(91, '08:55:25', '09:16', '16_SANTA ANA', 21, 'LP1912'),
    @(92, '08:55:25', '09:17', '27_EL RETIRO', 22, 'LP1912'),
    @(93, '08:31:53', '09:20', '26_HERNANDEZ', 49, 'LP1912'),
    @(94, '07:24:45', '09:21', '26_HERNANDEZ', 117, 'LP1912'),
    @(95, '08:00:32', '09:22', '16_SANTA ANA', 82, 'LP1912'),
    @(96, '07:24:45', '09:22', '17_ROMERO', 118, 'LP1912'),
    @(97, '07:48:31', '09:23', '11_ETCHEVERRY', 95, 'LP1912'),
    @(98, '07:48:31', '09:32', '15_ABASTO', 104, 'LP1912'),
    @(99, '07:48:31', '09:33', '10_OLMOS', 105, 'LP1912'),
    @(100, '08:47:51', '09:34', '23_HERNANDEZ', 47, 'LP1912'),
    @(101, '08:31:53', '09:41', '215C_EL PATO', 70, 'LP1912'),
    @(102, '07:48:31', '09:42', '215C_EL PATO', 114, 'LP1912'),
    @(103, '08:00:32', '09:43', '14_ABASTO', 103, 'LP1912'),
    @(104, '08:31:53', '09:46', '16_SANTA ANA', 75, 'LP1912'),
    @(105, '08:47:51', '09:52', '15_ABASTO', 65, 'LP1912'),
    @(106, '08:47:51', '09:53', '10_OLMOS', 66, 'LP1912'),
    @(107, '09:26:30', '09:58', '16_SANTA ANA', 32, 'LP1912'),
    @(108, '08:31:53', '10:03', '11_ETCHEVERRY', 92, 'LP1912'),
    @(109, '09:26:30', '10:04', '23_HERNANDEZ', 38, 'LP1912'),
    @(110, '08:31:53', '10:10', '16_P MOR-SANTA ANA', 99, 'LP1912'),
    @(111, '08:31:53', '10:12', '15_ABASTO', 101, 'LP1912'),
    @(112, '09:26:30', '10:13', '10_OLMOS', 47, 'LP1912'),
    @(113, '08:31:53', '10:20', '26_HERNANDEZ', 109, 'LP1912'),
    @(114, '08:55:25', '10:21', '26_HERNANDEZ', 86, 'LP1912'),
    @(115, '10:13:53', '10:22', '16_SANTA ANA', 9, 'LP1912'),
    @(116, '08:31:53', '10:22', '17_ROMERO', 111, 'LP1912'),
    @(117, '09:26:30', '10:23', '11_ETCHEVERRY', 57, 'LP1912'),
    @(118, '08:31:53', '10:26', '215A_EL PATO', 115, 'LP1912'),
    @(119, '10:13:53', '10:32', '10_OLMOS', 19, 'LP1912'),
    @(120, '10:13:53', '10:34', '16_SANTA ANA', 21, 'LP1912'),
    @(121, '10:13:53', '10:34', '23_HERNANDEZ', 21, 'LP1912'),
    @(122, '08:47:51', '10:41', '17_ROMERO', 114, 'LP1912'),
    @(123, '08:55:25', '10:42', '17_ROMERO', 107, 'LP1912'),
    @(124, '08:47:51', '10:43', '14_ABASTO', 116, 'LP1912'),
    @(125, '10:52:37', '10:52', '16_SANTA ANA', 0, 'LP1912'),
    @(126, '10:13:53', '10:52', '15_ABASTO', 39, 'LP1912'),
    @(127, '10:52:37', '10:53', '16_SANTA ANA', 1, 'LP1912'),
    @(128, '10:13:53', '10:56', '27_EL RETIRO', 43, 'LP1912'),
    @(129, '09:26:30', '10:57', '27_EL RETIRO', 91, 'LP1912'),
    @(130, '10:13:53', '11:01', '215C_EL PATO', 48, 'LP1912'),
    @(131, '09:26:30', '11:02', '215C_EL PATO', 96, 'LP1912'),
    @(132, '10:13:53', '11:03', '11_ETCHEVERRY', 50, 'LP1912'),
    @(133, '09:26:30', '11:06', '16_P MOR-167 Y 521', 100, 'LP1912'),
    @(134, '10:52:37', '11:06', '23_HERNANDEZ', 14, 'LP1912'),
    @(135, '10:13:53', '11:10', '10_OLMOS', 57, 'LP1912'),
    @(136, '10:52:37', '11:11', '10_OLMOS', 19, 'LP1912'),
    @(137, '10:13:53', '11:12', '15_ABASTO', 59, 'LP1912'),
    @(138, '10:13:53', '11:12', '23_HERNANDEZ', 59, 'LP1912'),
    @(139, '09:26:30', '11:19', '86_EST CHICA-ESC AGRARIA', 113, 'LP1912'),
    @(140, '10:13:53', '11:20', '26_HERNANDEZ', 67, 'LP1912'),
    @(141, '09:26:30', '11:21', '26_HERNANDEZ', 115, 'LP1912'),
    @(142, '10:13:53', '11:26', '225_C ROCA-H SUR', 73, 'LP1912'),
    @(143, '10:52:37', '11:27', '225_C ROCA-H SUR', 35, 'LP1912'),
    @(144, '10:13:53', '11:32', '81_EL PELIGRO', 79, 'LP1912'),
    @(145, '10:52:37', '11:34', '23_HERNANDEZ', 42, 'LP1912'),
    @(146, '10:52:37', '11:35', '11_ETCHEVERRY', 43, 'LP1912'),
    @(147, '10:13:53', '11:38', '10_OLMOS', 85, 'LP1912'),
    @(148, '10:13:53', '11:41', '17_ROMERO', 88, 'LP1912'),
    @(149, '10:52:37', '11:42', '17_ROMERO', 50, 'LP1912'),
    @(150, '10:52:37', '11:43', '10_OLMOS', 51, 'LP1912'),
    @(151, '10:13:53', '11:51', '215B_EL PATO', 98, 'LP1912'),
    @(152, '10:13:53', '11:58', '225_GOMEZ', 105, 'LP1912'),
    @(153, '10:52:37', '11:59', '225_GOMEZ', 67, 'LP1912'),
    @(154, '10:13:53', '12:02', '84_COLONIA URQUIZA-ESC 49', 109, 'LP1912'),
    @(155, '10:52:37', '12:06', '10_OLMOS', 74, 'LP1912'),
    @(156, '10:13:53', '12:06', '16_P MOR-SANTA ANA', 113, 'LP1912'),
    @(157, '10:13:53', '12:06', '14_ABASTO', 113, 'LP1912'),
    @(158, '10:52:37', '12:16', '14_ABASTO', 84, 'LP1912'),
    @(159, '10:52:37', '12:20', '215A_EL PATO', 88, 'LP1912'),
    @(160, '10:52:37', '12:21', '26_HERNANDEZ', 89, 'LP1912'),
    @(161, '10:52:37', '12:23', '17_ROMERO', 91, 'LP1912'),
    @(162, '10:52:37', '12:36', '27_EL RETIRO', 104, 'LP1912'),
    @(163, '10:52:37', '12:38', '17_179 Y 38', 106, 'LP1912'),
    @(164, '10:52:37', '12:50', '15_ABASTO', 118, 'LP1912')
)

foreach ($row in $data1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
}

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 10:52:37"
$ws2.Range("A3").Value = "Total filas: 21"

$data2 = @(
    @(6, '00:09:19', '01:12', '215_ALUAR', 63, 'LP1912'),
    @(7, '02:17:56', '02:57', '215_ALUAR', 40, 'LP1912'),
    @(8, '01:16:09', '02:58', '215_ALUAR', 102, 'LP1912'),
    @(9, '04:41:47', '04:45', '215A_EL PATO', 4, 'LP1912'),
    @(10, '02:56:23', '04:46', '215A_EL PATO', 110, 'LP1912'),
    @(11, '04:41:47', '05:34', '215B_EL PATO', 53, 'LP1912'),
    @(12, '03:55:07', '05:35', '215B_EL PATO', 100, 'LP1912'),
    @(13, '04:41:47', '06:11', '215A_EL PATO', 90, 'LP1912'),
    @(14, '04:19:04', '06:12', '215A_EL PATO', 113, 'LP1912'),
    @(15, '04:54:25', '06:46', '215C_EL PATO', 112, 'LP1912'),
    @(16, '05:18:42', '07:11', '215A_EL PATO', 113, 'LP1912'),
    @(17, '05:53:46', '07:51', '215D_EL PATO', 118, 'LP1912'),
    @(18, '06:44:40', '08:23', '215B_EL PATO', 99, 'LP1912'),
    @(19, '07:24:45', '09:01', '215A_EL PATO', 97, 'LP1912'),
    @(20, '08:31:53', '09:41', '215C_EL PATO', 70, 'LP1912'),
    @(21, '07:48:31', '09:42', '215C_EL PATO', 114, 'LP1912'),
    @(22, '08:31:53', '10:26', '215A_EL PATO', 115, 'LP1912'),
    @(23, '10:13:53', '11:01', '215C_EL PATO', 48, 'LP1912'),
    @(24, '09:26:30', '11:02', '215C_EL PATO', 96, 'LP1912'),
    @(25, '10:13:53', '11:51', '215B_EL PATO', 98, 'LP1912'),
    @(26, '10:52:37', '12:20', '215A_EL PATO', 88, 'LP1912')
)

foreach ($row in $data2) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 2).Value = $row[2]
    $ws2.Cells.Item($r, 3).Value = $row[3]
    $ws2.Cells.Item($r, 4).Value = $row[4]
    $ws2.Cells.Item($r, 5).Value = $row[5]
}

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 10:52:37"

Write-Host "Edit complete"
